# Updates the "cryptos" price/volume table with the latest scraped values.
# Values in column D that look like plain numbers are written with a leading
# apostrophe so Excel keeps them as text (matching the original inlineStr
# cells) instead of silently reformatting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.333.41"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.935.70"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'0.7595"
$ws.Range("E5").Value = "  +6.17%  "
$ws.Range("D6").Value = "'244.97"
$ws.Range("E6").Value = "  -2.73%  "
$ws.Range("D7").Value = "'0.9996"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "'27.73"
$ws.Range("E8").Value = "  +1.37%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.3185"
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("D10").Value = "'0.06998"
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("D11").Value = "'0.7796"
$ws.Range("D12").Value = "'0.08005"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("D13").Value = "1.936.71"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "'5.352"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "'94.31"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "'14.41"
$ws.Range("E16").Value = "  -2.79%  "
$ws.Range("D17").Value = "30.327.90"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "'252.74"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").Value = "'0.000007932"
$ws.Range("E19").Value = "  -2.65%  "
$ws.Range("D20").Value = "'5.759"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").Value = "2.185.27"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "'0.9994"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "'0.9980"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").Value = "'6.669"
$ws.Range("E24").Value = "  -3.52%  "
$ws.Range("D25").Value = "'9.470"
$ws.Range("E25").Value = "  -2.57%  "
$ws.Range("D26").Value = "'165.49"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").Value = "'18.98"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").Value = "'0.1332"
$ws.Range("E28").Value = "  +4.02%  "
$ws.Range("D29").Value = "'2.191"
$ws.Range("E29").Value = "  -5.33%  "
$ws.Range("D30").Value = "'1.365"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("D32").Value = "'4.396"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("D33").Value = "'4.125"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").Value = "'0.05154"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").Value = "'1.285"
$ws.Range("E35").Value = "  +1.30%  "
$ws.Range("D36").Value = "'0.7509"
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("D37").Value = "'2.768"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").Value = "'2.797"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "'77.55"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("D41").Value = "'6.410"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").Value = "'0.4457"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("D43").Value = "'1.966"
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("D44").Value = "'0.9997"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "'0.8344"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D46").Value = "'100.70"
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("D47").Value = "'9.755"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "'7.477"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").Value = "'37.43"
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("D50").Value = "'980.01"
$ws.Range("E50").Value = "  +11.13%  "
$ws.Range("D51").Value = "'0.06008"
$ws.Range("E51").Value = "  -1.40%  "
